# BetaFiberA-HW03.xlsx — "Updated notebook, reran simulation"
#
# Two new simulation configurations ("Holden" and "Rizzie Spiral") were
# inserted into the table right after "Spiral5" (i.e. as new rows 4 and 5),
# pushing all the subsequent rows down by two. The re-run also produced new
# numbers for the two brand new rows, while every pre-existing row kept its
# own label (column A index + column B name) but travelled down with the
# insert, carrying its own original numbers along with it. Separately,
# "Thomas Hex" was renamed to "Matthies Hex".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank rows at position 4 (old row 4 "RotRing OmegaMax-90"
#    and everything below it shifts down by two, to rows 6..31).
$ws.Rows("4:5").Insert()

# Excel's default "format like row above" for a 2-row insert did not bring
# the bold/centered/bordered style used by column A (style used by A2:A29)
# all the way through, so make sure A4:A5 match the rest of column A.
$ws.Range("A4:A5").Font.Bold = $true
$ws.Range("A4:A5").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A4:A5").VerticalAlignment = -4160     # xlTop
$ws.Range("A4:A5").Borders.LineStyle = 1

# 2) Fill in the new row 4 ("Holden") and row 5 ("Rizzie Spiral").
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$labels = New-Object 'object[,]' 2,1
$labels[0,0] = "Holden"
$labels[1,0] = "Rizzie Spiral"
$ws.Range("B4:B5").Value = $labels

$row4 = New-Object 'object[,]' 1,21
$row4[0,0]  = 0.8504883627377234
$row4[0,1]  = 0.9202821992972999
$row4[0,2]  = 1.766218421295849
$row4[0,3]  = 0.8504883627377234
$row4[0,4]  = 0.8828649458749206
$row4[0,5]  = 0.8828649458749117
$row4[0,6]  = 0.8828649458749117
$row4[0,7]  = 0.9455673715331205
$row4[0,8]  = 1.058153971255633
$row4[0,9]  = 1.408688568020527
$row4[0,10] = 0.8666190068804516
$row4[0,11] = 0.8828649458749117
$row4[0,12] = 1.766218421295849
$row4[0,13] = 1.308353392016786
$row4[0,14] = 1.355892896414485
$row4[0,15] = 1.166523909969495
$row4[0,16] = 1.187424718522231
$row4[0,17] = 1.166523909969495
$row4[0,18] = 1.111284775360401
$row4[0,19] = 1.065600809463303
$row4[0,20] = 1.087360355861939
$ws.Range("C4:W4").Value = $row4

$row5 = New-Object 'object[,]' 1,21
$row5[0,0]  = 1.583365550411374
$row5[0,1]  = 0.4384645455545183
$row5[0,2]  = 3.032117886040465
$row5[0,3]  = 1.583365550411374
$row5[0,4]  = 3.252219809860515
$row5[0,5]  = 3.252219809860515
$row5[0,6]  = 3.252219809860515
$row5[0,7]  = 1.097310259366113
$row5[0,8]  = 2.10178377575004
$row5[0,9]  = 2.188563281761665
$row5[0,10] = 0.997834488934809
$row5[0,11] = 3.252219809860515
$row5[0,12] = 3.032117886040465
$row5[0,13] = 2.307741718225919
$row5[0,14] = 2.064714072703289
$row5[0,15] = 2.622567748770785
$row5[0,16] = 1.90426456527265
$row5[0,17] = 2.622567748770785
$row5[0,18] = 2.241253376419616
$row5[0,19] = 2.443446663107796
$row5[0,20] = 1.836457449709938
$ws.Range("C5:W5").Value = $row5

# 3) Rename "Thomas Hex" -> "Matthies Hex" (now a few rows further down,
#    having shifted with the rest of the table).
$found = $ws.Cells.Find("Thomas Hex")
if ($found) {
    $found.Value = "Matthies Hex"
}
